$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "C"
$ws.Range("C22").Value = "C"
$ws.Range("D22").Value = "used left outer join fyi"
$ws.Range("C33").Value = "C"
$ws.Range("C34").Value = "C"
$ws.Range("C46").Value = "C"
